$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.527.05'
$ws.Range("E2").Value = '  -0.46%  '

$ws.Range("D3").Value = '1.673.06'
$ws.Range("E3").Value = '  -0.46%  '

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.010'
$c.Style = "Normal"
$ws.Range("E4").Value = '  +0.75%  '

$ws.Range("E5").Value = '  +0.82%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '307.38'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.52%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.3681'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -0.06%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '47.92'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -0.43%  '

$ws.Range("E9").Value = '  -1.52%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '1.174'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +1.28%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.07310'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +1.35%  '

$ws.Range("E12").Value = '  +0.64%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '6.165'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +1.40%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '20.46'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +1.90%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '6.788'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +1.59%  '

$ws.Range("D16").Value = '1.671.40'
$ws.Range("E16").Value = '  -0.52%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.00001094'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -0.50%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.06644'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.03%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '1.003'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +0.78%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '81.38'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +1.05%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '16.76'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +2.62%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '6.198'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +2.34%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '12.68'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +5.14%  '

$ws.Range("D24").Value = '24.550.53'
$ws.Range("E24").Value = '  -0.17%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.434'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +1.35%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '2.683'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +1.30%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '19.78'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +2.07%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '148.50'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -2.50%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '129.83'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +2.14%  '

$ws.Range("D30").Value = '1.859.31'
$ws.Range("E30").Value = '  -0.40%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.197'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +23.03%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '6.481'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +4.35%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '4.147'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +3.15%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.08594'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +2.39%  '

$ws.Range("B35").Value = 'Aptos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '13.20'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +7.55%  '

$ws.Range("B36").Value = 'WEMIXTOKEN'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '1.716'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +1.96%  '

$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.06469'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +1.60%  '

$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '5.392'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +2.08%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '8.839'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +2.04%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.02332'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +1.61%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.2161'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +3.85%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '1.232'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +0.01%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.6231'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +2.77%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '1.004'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +0.85%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '13.43'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +3.62%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '3.780'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +0.78%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.5910'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +1.13%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '2.042'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +2.29%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '125.75'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +0.29%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.07122'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -1.28%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '76.71'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +1.69%  '
